$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45889, 4535.36362203551, 3950.00607634303, 4224, 4904.608208, 3.96877759614631),
    @(45890, 4535.36362203551, 3958.32134296011, 1944, 4904.608208, 99.3152470385248),
    @(45891, 4535.36362203551, 3959.84488557803, 1944, 4904.608208, 99.3787279809382),
    @(45892, 868.670076800243, 1547.20039187832, 1944, 1707.033454, 18.3984903782534),
    @(45893, 773.003570643336, 1535.23478929175, 1944, 1602.771126, 17.5417643603504),
    @(45894, 5476.20428821302, 4288.16179237205, 1944, 5817.596568, 111.89808633996),
    @(45895, 5476.20428821302, 4311.74729772894, 1944, 5817.596568, 112.88081572983),
    @(45896, 5476.20428821302, 4231.98262777828, 1944, 5817.596568, 109.557287815219),
    @(45897, 5476.20428821302, 4211.24853044701, 1944, 5817.596568, 108.693367093083),
    @(45898, 5476.20428821302, 4242.23073439478, 1944, 5817.596568, 109.984292257574),
    @(45899, 948.243530910975, 1634.87752733635, 1944, 1809.959712, 23.0247378510575),
    @(45900, 832.128236149678, 1598.63493809623, 1944, 1682.357187, 21.0359953727731),
    @(45901, 5776.55841866516, 4989.77305918669, 2952, 6390.755333, 110.49874889673),
    @(45902, 5784.11073458408, 4979.27710125451, 2952, 6398.371705, 110.064086319601)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
